# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Poroto verde" (Agricola del Norte S.A.
# de Arica) at row 31, shifting the existing rows 31-42 down to 32-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31, pushing rows 31..42 down to 32..43.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new data point.
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 44529
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 100112031
$ws.Range("G31").Value = "Poroto verde"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 1700
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 650
$ws.Range("M31").Value = 625
$ws.Range("N31").Value = "`$/kilo"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 625
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
